$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "gpt-4o-2024-05-13",
    0.96,
    0.96,
    1,
    0.9399999999999999,
    0.86,
    0.98,
    0.95,
    0.87,
    0.92,
    0.895,
    0.68,
    0.6,
    0.64,
    0.8,
    0.76,
    0.861,
    0.82,
    0.9399999999999999,
    0.96,
    0.907,
    0.367,
    0.6830000000000001,
    0.4,
    0.57,
    0.384,
    0.8
)

$row = 3
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($row, $col).Value = $values[$i]
}

# Match the formatting (style) used by the row above it.
$ws.Range("A2:AA2").Copy() | Out-Null
$ws.Range("A3:AA3").PasteSpecial(-4122) | Out-Null
